$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DiemTB (N) and DiemMonNV (W) for row 3 and row 4,
# and DiemMon3 (Y) for row 4 - switching it from a number to a text value.
$ws.Range("N3").Value = "7,8"
$ws.Range("W3").Value = "8,9"
$ws.Range("N4").Value = "7,9"
$ws.Range("Y4").Value = "6,6"

# Update the view/selection state of the sheet.
$ws.Application.ActiveWindow.ScrollColumn = 17
$ws.Range("AA7").Select()
